$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "magapoke_2026-02-25"

$data = @(
    @('rank', 'title'),
    @(1, 'アイドラトリィ'),
    @(2, 'ペンの夢に紅をさす'),
    @(3, '黒月のイェルクナハト'),
    @(4, 'あの島の海音荘'),
    @(5, 'ドリーム☆ジャンボ☆ガール'),
    @(6, '黄昏町プリズナーズ'),
    @(7, 'K-9~警視庁公安部公安第9課異能対策係~'),
    @(8, '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜'),
    @(9, 'ゼロとヒャク'),
    @(10, 'ハードワーカー中田'),
    @(11, 'せいぶつ部の田辺くん'),
    @(12, 'ハンドレッドノート－アグリーダック－'),
    @(13, 'ともだちづくり'),
    @(14, '怨霊日和'),
    @(15, '普通の本はありません！'),
    @(16, 'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'),
    @(17, '白鳥運子は31画'),
    @(18, 'お前がヤったんだろ！'),
    @(19, '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'),
    @(20, '屋根の下のアルテミス'),
    @(21, 'MYS'),
    @(22, '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'),
    @(23, 'わが投資術　市場は誰に微笑むか'),
    @(24, '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'),
    @(25, '君が監督！'),
    @(26, 'その青春'),
    @(27, '鳴るさんだぁ'),
    @(28, 'ナマイキ旭ちゃんをわからせたい'),
    @(29, '明智ナンバーワン'),
    @(30, '春くらり'),
    @(31, '平成転生'),
    @(32, '追放されなかった男　～二度目の人生は土下座から始まりました～'),
    @(33, 'じゅーくぼっくす'),
    @(34, '歪みの虜'),
    @(35, '人生逆転ダンジョン'),
    @(36, 'JK Biker'),
    @(37, 'GURU'),
    @(38, 'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'),
    @(39, 'ch登録お願いします！'),
    @(40, '眠れる森のレガ'),
    @(41, 'ハプスブルク家の華麗なる受難'),
    @(42, '邪目さんは邪神です'),
    @(43, '白銀のキュイジーヌ～明治外交官の料理人～'),
    @(44, '永久のユウグレ'),
    @(45, '宇曽田みのりの代用料理'),
    @(46, '花子狩り'),
    @(47, 'イエティ、とある日々'),
    @(48, 'きゃわるり方程式')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}